# Applies the 2024-09-18 04:59:21 update to the "2024" sheet:
#  1. A new September transaction note ("your relationship") is logged at
#     10:29:06, pushing the existing September_Details/September_Date
#     history (previously rows 43-156) down by one row (to rows 44-157).
#  2. The trailing "Broadband" group label moves from A165 to a new A166,
#     growing the used range from A1:Y165 to A1:Y166.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# --- September (columns R="September_Details", S="September_Date") ---
# Each tuple is (row, details, date) for the post-edit state, newest first.
$septData = @(
    ,@(43, "your relationship", "2024-09-18 10:29:06")
    ,@(44, "balance your axis", "2024-09-18 10:28:28")
    ,@(45, "axis", "2024-09-18 08:12:44")
    ,@(46, "broker", "2024-09-18 04:09:58")
    ,@(47, "balance your axis", "2024-09-17 13:07:16")
    ,@(48, "dispute", "2024-09-16 12:53:44")
    ,@(49, "money google icici", "2024-09-16 12:53:29")
    ,@(50, "indusind", "2024-09-16 11:13:15")
    ,@(51, "communication feedback", "2024-09-16 11:13:15")
    ,@(52, "balance your axis", "2024-09-16 08:57:11")
    ,@(53, "balance your axis", "2024-09-16 07:57:00")
    ,@(54, "money google icici", "2024-09-15 21:06:00")
    ,@(55, "adani icici", "2024-09-15 13:10:50")
    ,@(56, "balance your axis", "2024-09-15 07:56:24")
    ,@(57, "bal axisbank w axis", "2024-09-15 07:12:01")
    ,@(58, "hdfc", "2024-09-14 21:25:23")
    ,@(59, "change the", "2024-09-12 21:16:38")
    ,@(60, "dispute", "2024-09-12 19:02:14")
    ,@(61, "congrats limit icici", "2024-09-12 19:03:39")
    ,@(62, "latest transaction pan", "2024-09-12 12:22:12")
    ,@(63, "assistance", "2024-09-12 12:17:33")
    ,@(64, "balance your axis", "2024-09-12 09:37:28")
    ,@(65, "bal axisbank", "2024-09-12 00:54:39")
    ,@(66, "your relationship", "2024-09-11 16:05:27")
    ,@(67, "login internet personal share", "2024-09-11 14:16:45")
    ,@(68, "balance your axis", "2024-09-11 12:45:33")
    ,@(69, "balance your axis", "2024-09-11 09:45:01")
    ,@(70, "axis", "2024-09-11 06:57:42")
    ,@(71, "money google icici", "2024-09-10 20:42:12")
    ,@(72, "dispute", "2024-09-10 20:42:34")
    ,@(73, "reward points cash", "2024-09-10 19:43:35")
    ,@(74, "balance your axis", "2024-09-10 13:32:42")
    ,@(75, "ach indianesign bal axisbank", "2024-09-10 13:22:37")
    ,@(76, "ach indianesign bal axisbank", "2024-09-10 13:22:37")
    ,@(77, "balance your axis", "2024-09-10 11:21:40")
    ,@(78, "your relationship", "2024-09-10 11:02:23")
    ,@(79, "bank bal broker", "2024-09-09 19:59:02")
    ,@(80, "beneficiary", "2024-09-09 15:48:10")
    ,@(81, "beneficiary saravanan", "2024-09-09 14:52:20")
    ,@(82, "bal axisbank", "2024-09-09 12:19:34")
    ,@(83, "bal axisbank", "2024-09-09 12:19:33")
    ,@(84, "dispute", "2024-09-09 12:17:30")
    ,@(85, "bal axisbank", "2024-09-09 12:04:31")
    ,@(86, "transfer freedom share anyone axis", "2024-09-09 11:56:19")
    ,@(87, "corporate internet share", "2024-09-09 11:40:49")
    ,@(88, "corporate internet share", "2024-09-09 11:39:30")
    ,@(89, "bal axisbank", "2024-09-09 11:38:16")
    ,@(90, "bal axisbank", "2024-09-09 11:38:16")
    ,@(91, "bal axisbank", "2024-09-09 11:38:15")
    ,@(92, "bal axisbank", "2024-09-09 11:38:15")
    ,@(93, "corporate internet share", "2024-09-09 11:35:34")
    ,@(94, "corporate internet share", "2024-09-09 11:32:23")
    ,@(95, "ift anbu tpar", "2024-09-09 11:27:52")
    ,@(96, "balance your axis", "2024-09-09 11:24:00")
    ,@(97, "corporate internet share", "2024-09-09 11:21:43")
    ,@(98, "corporate internet share", "2024-09-09 11:17:34")
    ,@(99, "corporate internet share", "2024-09-09 11:15:51")
    ,@(100, "corporate internet share", "2024-09-09 11:14:13")
    ,@(101, "anbu tparty bal axisbank", "2024-09-09 11:13:37")
    ,@(102, "corporate internet share", "2024-09-09 11:10:39")
    ,@(103, "corporate internet share", "2024-09-09 11:07:31")
    ,@(104, "corporate internet share", "2024-09-09 11:03:09")
    ,@(105, "saravanan", "2024-09-09 10:43:11")
    ,@(106, "balance your axis", "2024-09-09 08:10:16")
    ,@(107, "ekalaivan", "2024-09-08 18:40:34")
    ,@(108, "balance your axis", "2024-09-08 09:53:37")
    ,@(109, "balance your axis", "2024-09-07 12:12:22")
    ,@(110, "balance your axis", "2024-09-07 09:34:58")
    ,@(111, "bal axis", "2024-09-07 08:46:40")
    ,@(112, "axis", "2024-09-07 08:31:28")
    ,@(113, "your relationship", "2024-09-06 12:23:25")
    ,@(114, "balance your axis", "2024-09-06 09:55:31")
    ,@(115, "beneficiary", "2024-09-05 17:13:56")
    ,@(116, "coimbatore ramalinga", "2024-09-05 17:06:01")
    ,@(117, "beneficiary", "2024-09-05 17:04:10")
    ,@(118, "bal axisbank", "2024-09-05 16:52:25")
    ,@(119, "share anyone axis", "2024-09-05 16:38:59")
    ,@(120, "transfer anyone axis", "2024-09-05 16:35:58")
    ,@(121, "share anyone axis", "2024-09-05 16:31:34")
    ,@(122, "transfer", "2024-09-05 16:28:38")
    ,@(123, "bal axisbank axis", "2024-09-05 16:26:56")
    ,@(124, "bal axisbank", "2024-09-05 16:26:55")
    ,@(125, "transfer", "2024-09-05 16:25:07")
    ,@(126, "transfer", "2024-09-05 16:22:23")
    ,@(127, "share anyone axis", "2024-09-05 16:06:05")
    ,@(128, "internet bal axisbank", "2024-09-05 16:05:55")
    ,@(129, "transfer share anyone axis", "2024-09-05 16:03:14")
    ,@(130, "axis", "2024-09-05 15:57:15")
    ,@(131, "your net internet", "2024-09-05 15:57:15")
    ,@(132, "hear your feedback atm", "2024-09-05 14:21:08")
    ,@(133, "axis bna", "2024-09-05 14:18:32")
    ,@(134, "axis bna", "2024-09-05 14:13:16")
    ,@(135, "axis bna", "2024-09-05 14:15:23")
    ,@(136, "balance your axis", "2024-09-05 09:20:57")
    ,@(137, "bal axis", "2024-09-05 09:06:25")
    ,@(138, "broker", "2024-09-04 21:20:47")
    ,@(139, "exclusive on axis", "2024-09-04 13:21:05")
    ,@(140, "your corporate axis", "2024-09-04 11:46:10")
    ,@(141, "balance your axis", "2024-09-04 08:14:16")
    ,@(142, "axis", "2024-09-04 07:02:13")
    ,@(143, "bal axisbank w axis", "2024-09-04 06:53:15")
    ,@(144, "logging iob internet", "2024-09-03 20:09:12")
    ,@(145, "password internet", "2024-09-03 20:05:31")
    ,@(146, "logging iob internet", "2024-09-03 20:05:09")
    ,@(147, "internet", "2024-09-03 19:58:18")
    ,@(148, "login internet invalid", "2024-09-03 19:54:49")
    ,@(149, "login internet invalid", "2024-09-03 19:56:17")
    ,@(150, "corporate internet share", "2024-09-03 19:22:58")
    ,@(151, "login sbi internet personal do not share anyone", "2024-09-03 19:17:10")
    ,@(152, "login internet personal share", "2024-09-03 19:13:40")
    ,@(153, "internet verify it", "2024-09-03 19:05:49")
    ,@(154, "balance your axis", "2024-09-03 13:14:06")
    ,@(155, "lounge", "2024-09-03 13:08:08")
    ,@(156, "balance your axis", "2024-09-03 11:21:30")
    ,@(157, "broker", "2024-09-01 22:35:38")
)

foreach ($entry in $septData) {
    $rowNum = $entry[0]
    $ws.Range("R$rowNum").Value = $entry[1]
    $ws.Range("S$rowNum").Value = $entry[2]
}

# --- "Broadband" group label shifts from row 165 to the new row 166 ---
$ws.Range("A165").Value = ""
$ws.Range("A166").Value = "Broadband"

